# "connected excel for valid credential test"
#
# - Remove the unused blank "Sheet2" tab entirely, and have the data that
#   used to live in "validCredentialTest" be replaced by a freshly
#   populated sheet (re-using the blank "Sheet2" tab, renamed back to
#   "validCredentialTest") that now contains the Admin test-data table.
# - "Sheet3" (the other trailing blank sheet) is left alone.

$wb = $excel.ActiveWorkbook

# Drop the old (empty) "validCredentialTest" sheet and reuse the blank
# "Sheet2" tab in its place, renamed back to "validCredentialTest" - this
# keeps the original sheet order (invalidCredentialTest, validCredentialTest,
# Sheet3) while dropping the extra blank tab.
$old = $wb.Worksheets.Item("validCredentialTest")
[void]$old.Delete()

$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "validCredentialTest"

# Populate the valid-credential test data table.
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Expected Value"

$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Range("C2").Value = "Employee List"

$ws.Range("A3").Value = "Admin"
$ws.Range("B3").Value = "admin123"
$ws.Range("C3").Value = "Employee List"

# Widen column C to fit its longest value and leave the last data row
# selected, same as the authored workbook.
$ws.Columns.Item(3).AutoFit()
[void]$ws.Range("A3:C3").Select()
